# "Added CDS All studies testcase"
# The "Sample ID" SQL query stored in cell B3 is trimmed: the Tumor and
# Analyte Type columns are removed from its SELECT clause. The active
# selection also moves from C4 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
 s.phs_accession = 'phs001437' AND gi.library_selection = 'cDNA'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# Strip the single trailing newline the here-string literal adds after the
# last line so the stored text matches the original cell exactly.
$newQuery = $newQuery -replace "`r?`n$", ""

$ws.Range("B3").Value = $newQuery

$ws.Range("C3").Select()
